$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (Through 2022-07-15 -> Through 2022-07-16)
$ws.Name = "Through 2022-07-16"

# Update header label for July
$ws.Range("A8").Value = "July (through 07-16)"

# Row 8 (July) updated values for columns D:I (2017-2022)
$ws.Range("D8").Value = 33
$ws.Range("E8").Value = 38
$ws.Range("F8").Value = 26
$ws.Range("G8").Value = 65
$ws.Range("H8").Value = 75
$ws.Range("I8").Value = 91

# Row 9 (Total) updated values for columns D:I (2017-2022)
$ws.Range("D9").Value = 423
$ws.Range("E9").Value = 391
$ws.Range("F9").Value = 277
$ws.Range("G9").Value = 537
$ws.Range("H9").Value = 835
$ws.Range("I9").Value = 896
